# Hotels/review sheets were swapped: the tab that used to hold the
# hotel summary ("hotel_info", 1st tab) now becomes the review-level
# sheet ("review_info"), and the tab that used to hold the (empty)
# review headers ("review_info", 2nd tab) becomes the hotel summary
# sheet ("hotel_info"), gaining a new "State" column and the hotel data
# row that used to live on the first tab.

$wb = $excel.ActiveWorkbook

$s1 = $wb.Sheets.Item(1)
$s2 = $wb.Sheets.Item(2)

# --- rename tabs (swap names), via a scratch name to dodge collisions ---
$s1.Name = "tmp_swap_name"
$s2.Name = "hotel_info"
$s1.Name = "review_info"

# --- rebuild the review_info sheet (1st tab): header-only, 25 columns ---
$s1.Cells.Clear()
$reviewHeaders = @( `
    "STR", `
    "reviewer_ID", `
    "reviewer_name", `
    "Review_ID", `
    "Date_of_scraping", `
    "ReviewURL", `
    "Tripadvisor_gcode", `
    "Tripadvisor_dcode", `
    "Tripadvisor_rcode", `
    "review_date", `
    "review_title", `
    "review_content", `
    "review_rating", `
    "trip_month", `
    "trip_purpose", `
    "value", `
    "rooms", `
    "Location", `
    "Cleanliness", `
    "Sleep Quality", `
    "Service", `
    "Picture(yes=1)", `
    "respondent", `
    "response_date", `
    "response_text" `
)
for ($i = 0; $i -lt $reviewHeaders.Length; $i++) {
    $s1.Cells.Item(1, $i + 1).Value = $reviewHeaders[$i]
}

# --- rebuild the hotel_info sheet (2nd tab): header + one data row, with new "State" column ---
$s2.Cells.Clear()
$hotelHeaders = @("STR", "Hotel_Name", "State", "City", "Zip", "TA_ReviewURL", "Tripadvisor_Hotel_Name", "English_Reviews_num", "Local_Rank", "Total_Reviews_num")
for ($i = 0; $i -lt $hotelHeaders.Length; $i++) {
    $s2.Cells.Item(1, $i + 1).Value = $hotelHeaders[$i]
}

$s2.Cells.Item(2, 1).Value = 9793
$s2.Cells.Item(2, 2).Value = "Hilton New Orleans Riverside"
$s2.Cells.Item(2, 3).Value = "Louisiana"
$s2.Cells.Item(2, 4).Value = "New Orleans"
$s2.Cells.Item(2, 5).Value = 70130
$s2.Cells.Item(2, 6).Value = "https://www.tripadvisor.com/Hotel_Review-g60864-d93164-Reviews-Hilton_New_Orleans_Riverside-New_Orleans_Louisiana.html"
$s2.Cells.Item(2, 7).Value = "Hilton New Orleans Riverside"

# These three look numeric but must be stored as TEXT (matches the
# source data). Force text entry via NumberFormat, then restore the
# default "Normal" style so no stray per-cell formatting lingers.
foreach ($cell in @(
        @{ r = 2; c = 8; v = "5852" }, `
        @{ r = 2; c = 9; v = "86" }, `
        @{ r = 2; c = 10; v = "6018" } `
    )) {
    $rng = $s2.Cells.Item($cell.r, $cell.c)
    $rng.NumberFormat = "@"
    $rng.Value = $cell.v
    $rng.Style = "Normal"
}
